# Scheduled runner: refresh market-board derived columns (H:N) on a handful
# of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with newly
# polled prices. Row numbers below are the worksheet row numbers (the sheet
# header occupies row 1).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3810.0557
$ws.Range("I74").Value = 3300
$ws.Range("J74").Value = 3840.0588
$ws.Range("K74").Value = 3300
$ws.Range("L74").Value = 3840.0588
$ws.Range("M74").Value = -2364
$ws.Range("N74").Value = -5712.0588

$ws.Range("H77").Value = 3810.0557
$ws.Range("I77").Value = 3300
$ws.Range("J77").Value = 3840.0588
$ws.Range("K77").Value = 16500
$ws.Range("L77").Value = 19200.294
$ws.Range("M77").Value = -11820
$ws.Range("N77").Value = -28560.294

$ws.Range("H82").Value = 266.66666
$ws.Range("I82").Value = 266.66666
$ws.Range("K82").Value = 799.9999799999999
$ws.Range("M82").Value = -393.9999799999999

$ws.Range("H85").Value = 266.66666
$ws.Range("I85").Value = 266.66666
$ws.Range("K85").Value = 799.9999799999999
$ws.Range("M85").Value = 604.0000200000001

$ws.Range("H95").Value = 46899.75
$ws.Range("J95").Value = 46899.75
$ws.Range("L95").Value = 46899.75
$ws.Range("N95").Value = -52391.75

$ws.Range("H135").Value = 906.7838
$ws.Range("I135").Value = 898.6389
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 8087.7501
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -5552.7501
$ws.Range("N135").Value = -15870

$ws.Range("H137").Value = 1594.1471
$ws.Range("I137").Value = 1461.3226
$ws.Range("J137").Value = 2966.6667
$ws.Range("K137").Value = 4383.9678
$ws.Range("L137").Value = 8900.000100000001
$ws.Range("M137").Value = -1833.9678
$ws.Range("N137").Value = -14000.0001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 17001.6
$ws.Range("I26").Value = 3750
$ws.Range("K26").Value = 3750
$ws.Range("M26").Value = -3420

$ws.Range("H32").Value = 4144.9355
$ws.Range("I32").Value = 4085.228
$ws.Range("J32").Value = 4825.6
$ws.Range("K32").Value = 4085.228
$ws.Range("L32").Value = 4825.6
$ws.Range("M32").Value = -3798.228
$ws.Range("N32").Value = -5399.6

$ws.Range("H61").Value = 2272
$ws.Range("I61").Value = 1076.5714
$ws.Range("J61").Value = 3666.6667
$ws.Range("K61").Value = 1076.5714
$ws.Range("L61").Value = 3666.6667
$ws.Range("M61").Value = -864.5714
$ws.Range("N61").Value = -4090.6667

$ws.Range("H132").Value = 8173.6113
$ws.Range("I132").Value = 18135.666
$ws.Range("J132").Value = 3192.5833
$ws.Range("K132").Value = 54406.99800000001
$ws.Range("L132").Value = 9577.749899999999
$ws.Range("M132").Value = -51876.99800000001
$ws.Range("N132").Value = -14637.7499

$ws.Range("H136").Value = 2272
$ws.Range("I136").Value = 1076.5714
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 3229.7142
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -679.7142000000003
$ws.Range("N136").Value = -16100.0001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 52073.45
$ws.Range("I134").Value = 78397.62
$ws.Range("J134").Value = 3185.7144
$ws.Range("K134").Value = 235192.86
$ws.Range("L134").Value = 9557.143199999999
$ws.Range("M134").Value = -232657.86
$ws.Range("N134").Value = -14627.1432

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1390.6
$ws.Range("I58").Value = 1490.2222
$ws.Range("J58").Value = 1309.091
$ws.Range("K58").Value = 1490.2222
$ws.Range("L58").Value = 1309.091
$ws.Range("M58").Value = -1287.2222
$ws.Range("N58").Value = -1715.091

$ws.Range("H132").Value = 2909.0454
$ws.Range("I132").Value = 2006
$ws.Range("J132").Value = 5979.4
$ws.Range("K132").Value = 6018
$ws.Range("L132").Value = 17938.2
$ws.Range("M132").Value = -3488
$ws.Range("N132").Value = -22998.2

$ws.Range("H136").Value = 1390.6
$ws.Range("I136").Value = 1490.2222
$ws.Range("J136").Value = 1309.091
$ws.Range("K136").Value = 4470.6666
$ws.Range("L136").Value = 3927.273
$ws.Range("M136").Value = -1920.6666
$ws.Range("N136").Value = -9027.272999999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 133.14285
$ws.Range("I47").Value = 133.14285
$ws.Range("K47").Value = 399.42855
$ws.Range("M47").Value = 31.57144999999997

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 128588.81
$ws.Range("I132").Value = 156878.69
$ws.Range("K132").Value = 470636.07
$ws.Range("M132").Value = -468106.07

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1158.8334
$ws.Range("J68").Value = 1250
$ws.Range("L68").Value = 1250
$ws.Range("N68").Value = -2748

$ws.Range("H71").Value = 1158.8334
$ws.Range("J71").Value = 1250
$ws.Range("L71").Value = 6250
$ws.Range("N71").Value = -13738

$ws.Range("H82").Value = 1085.2
$ws.Range("I82").Value = 956.5
$ws.Range("K82").Value = 956.5
$ws.Range("M82").Value = -595.5

$ws.Range("H85").Value = 1085.2
$ws.Range("I85").Value = 956.5
$ws.Range("K85").Value = 956.5
$ws.Range("M85").Value = 291.5

$ws.Range("H122").Value = 3610.2666
$ws.Range("I122").Value = 4309.143
$ws.Range("J122").Value = 2998.75
$ws.Range("K122").Value = 12927.429
$ws.Range("L122").Value = 8996.25
$ws.Range("M122").Value = -10477.429
$ws.Range("N122").Value = -13896.25

$ws.Range("H132").Value = 3232.2917
$ws.Range("I132").Value = 2957.3333
$ws.Range("J132").Value = 3507.25
$ws.Range("K132").Value = 8871.999899999999
$ws.Range("L132").Value = 10521.75
$ws.Range("M132").Value = -6341.999899999999
$ws.Range("N132").Value = -15581.75

$ws.Range("H136").Value = 4605.636
$ws.Range("I136").Value = 5432.75
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 16298.25
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -13748.25
$ws.Range("N136").Value = -12300

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2683
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 2719.6
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 2719.6
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -3967.6

$ws.Range("H65").Value = 2683
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 2719.6
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 13598
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -19838

# Price dropped to 0 (no longer profitable) - clear the now-stale loss cell
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H113").Value = 589.9231
$ws.Range("I113").Value = 786
$ws.Range("J113").Value = 467.375
$ws.Range("K113").Value = 2358
$ws.Range("L113").Value = 1402.125
$ws.Range("M113").Value = -188
$ws.Range("N113").Value = -5742.125

$ws.Range("H132").Value = 2267.2
$ws.Range("I132").Value = 1356.6666
$ws.Range("J132").Value = 3633
$ws.Range("K132").Value = 4069.9998
$ws.Range("L132").Value = 10899
$ws.Range("M132").Value = -1539.9998
$ws.Range("N132").Value = -15959

$ws.Range("H136").Value = 5503.3105
$ws.Range("I136").Value = 6141.5
$ws.Range("J136").Value = 2440
$ws.Range("K136").Value = 18424.5
$ws.Range("L136").Value = 7320
$ws.Range("M136").Value = -15874.5
$ws.Range("N136").Value = -12420
